$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Update Players sheet (game status / box score updates for 2026-01-17)
$ws.Range("G2").Value = "Final"
$ws.Range("G3").Value = "11:51 - 2nd Half"
$ws.Range("G4").Value = "11:02 - 1st Half"
$ws.Range("H4").Value = 0
$ws.Range("M4").Value = 1
$ws.Range("O4").Value = 7
$ws.Range("G5").Value = "11:51 - 2nd Half"
$ws.Range("H5").Value = 18
$ws.Range("O5").Value = 25
$ws.Range("G6").Value = "Final"
$ws.Range("H6").Value = 22
$ws.Range("G7").Value = "11:02 - 1st Half"
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 3
$ws.Range("O7").Value = 8
$ws.Range("G8").Value = "Final"
$ws.Range("G9").Value = "Final"
$ws.Range("O9").Value = 27
$ws.Range("G10").Value = "11:02 - 1st Half"
$ws.Range("G11").Value = "11:51 - 2nd Half"
$ws.Range("H11").Value = -2
$ws.Range("O11").Value = 18
$ws.Range("G12").Value = "11:02 - 1st Half"
$ws.Range("H12").Value = 6
$ws.Range("G13").Value = "11:51 - 2nd Half"
$ws.Range("H13").Value = 4
$ws.Range("K13").Value = 2
$ws.Range("O13").Value = 20
$ws.Range("G14").Value = "Final"
$ws.Range("G15").Value = "11:51 - 2nd Half"
$ws.Range("G16").Value = "11:51 - 2nd Half"
$ws.Range("O16").Value = 11
$ws.Range("G17").Value = "11:51 - 2nd Half"
$ws.Range("G18").Value = "11:02 - 1st Half"
$ws.Range("G19").Value = "Final"
$ws.Range("G20").Value = "11:02 - 1st Half"
$ws.Range("G21").Value = "11:51 - 2nd Half"
$ws.Range("J21").Value = 4
$ws.Range("O21").Value = 19
$ws.Range("G22").Value = "11:51 - 2nd Half"
$ws.Range("H22").Value = 11
$ws.Range("J22").Value = 3
$ws.Range("M22").Value = 3
$ws.Range("O22").Value = 19
$ws.Range("G23").Value = "Final"
$ws.Range("J23").Value = 7
$ws.Range("G24").Value = "11:51 - 2nd Half"
$ws.Range("O24").Value = 25
$ws.Range("G25").Value = "11:51 - 2nd Half"
$ws.Range("O25").Value = 26
$ws.Range("G26").Value = "Final"
$ws.Range("H26").Value = 10
$ws.Range("M26").Value = 4
$ws.Range("G27").Value = "11:02 - 1st Half"
$ws.Range("O27").Value = 7
$ws.Range("G28").Value = "Final"
$ws.Range("G29").Value = "11:02 - 1st Half"
$ws.Range("H29").Value = 4
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = 2
$ws.Range("O29").Value = 6
$ws.Range("G30").Value = "11:51 - 2nd Half"
$ws.Range("O30").Value = 17
$ws.Range("G31").Value = "11:02 - 1st Half"
$ws.Range("H31").Value = 4
$ws.Range("J31").Value = 3
$ws.Range("O31").Value = 7
$ws.Range("G32").Value = "Final"
$ws.Range("G33").Value = "11:51 - 2nd Half"
$ws.Range("H33").Value = 24
$ws.Range("I33").Value = 19
$ws.Range("J33").Value = 7
$ws.Range("O33").Value = 22
$ws.Range("G34").Value = "11:51 - 2nd Half"
$ws.Range("G35").Value = "11:02 - 1st Half"
$ws.Range("H35").Value = 6
$ws.Range("O35").Value = 9
$ws.Range("G36").Value = "11:02 - 1st Half"
$ws.Range("G37").Value = "Final"
$ws.Range("G38").Value = "Final"
$ws.Range("G39").Value = "Final"
$ws.Range("G40").Value = "11:02 - 1st Half"
$ws.Range("O40").Value = 7
$ws.Range("D41").Value = "Isaiah Brown"
$ws.Range("E41").Value = "FLA"
$ws.Range("F41").Value = "FLA@VAN"
$ws.Range("G41").Value = "11:02 - 1st Half"
$ws.Range("H41").Value = 6
$ws.Range("K41").Value = 0
$ws.Range("O41").Value = 3
$ws.Range("D42").Value = "Andrija Jelavić"
$ws.Range("G42").Value = "Final"
$ws.Range("H42").Value = 5
$ws.Range("I42").Value = 5
$ws.Range("J42").Value = 1
$ws.Range("O42").Value = 12
$ws.Range("D43").Value = "Brandon Garrison"
$ws.Range("E43").Value = "UK"
$ws.Range("F43").Value = "UK@TENN"
$ws.Range("G43").Value = "Final"
$ws.Range("H43").Value = 4
$ws.Range("J43").Value = 2
$ws.Range("K43").Value = 1
$ws.Range("O43").Value = 18
$ws.Range("G44").Value = "Final"
$ws.Range("G45").Value = "Final"
$ws.Range("D46").Value = "London Jemison"
$ws.Range("E46").Value = "ALA"
$ws.Range("F46").Value = "ALA@OU"
$ws.Range("G46").Value = "11:51 - 2nd Half"
$ws.Range("J46").Value = 3
$ws.Range("M46").Value = 1
$ws.Range("O46").Value = 11
$ws.Range("D47").Value = "Tyler Harris"
$ws.Range("E47").Value = "VAN"
$ws.Range("G47").Value = "11:02 - 1st Half"
$ws.Range("I47").Value = 2
$ws.Range("K47").Value = 0
$ws.Range("G48").Value = "Final"
$ws.Range("D49").Value = "Noah Williamson"
$ws.Range("G49").Value = "11:51 - 2nd Half"
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("O49").Value = 3
$ws.Range("D50").Value = "Urban Klavzar"
$ws.Range("E50").Value = "FLA"
$ws.Range("F50").Value = "FLA@VAN"
$ws.Range("G50").Value = "11:02 - 1st Half"
$ws.Range("K50").Value = 1
$ws.Range("O50").Value = 5
$ws.Range("G51").Value = "Final"
$ws.Range("G52").Value = "Final"
$ws.Range("G53").Value = "Final"

# Update OwnerTotals sheet (starter_pooh_total changes)
$ws2.Range("B2").Value = 42
$ws2.Range("B3").Value = 33
$ws2.Range("B4").Value = 29
$ws2.Range("B5").Value = 24
$ws2.Range("B7").Value = 12
